# Updated cryptos list values (Price + Volume(1h)) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.122.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.66%  "
$ws.Range("D3").Value = "'2.432.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'317.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").Value = "'102.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.47%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +7.46%  "
$ws.Range("D10").Value = "'35.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("D13").Value = "'18.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("D14").Value = "'7.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "'2.814.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "'2.426.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "'45.050.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").Value = "'12.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "'6.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'0.0₃0923"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "'244.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").Value = "'2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "'49.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").Value = "'32.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "'20.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.44%  "
$ws.Range("D33").Value = "'0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.84%  "
$ws.Range("D34").Value = "'5.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.0765"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "'4.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("D40").Value = "'124.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.42%  "
$ws.Range("D41").Value = "'2.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "'20.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").Value = "'1.934.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "'2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").Value = "'9.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").Value = "'1.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.94%  "
$ws.Range("D50").Value = "'76.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.82%  "
$ws.Range("D51").Value = "'53.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.62%  "
